$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new TODO entries (column D) below the existing list ---
# Shared-string insertion order must match: Topic2Vec, LDA visualization from
# paper, LDA2Vec: does not work (import problem)  -> so set D68/D69 first,
# then D67, to get that exact shared-string ordering.
$ws.Range("D68").Value = "Topic2Vec"
$ws.Range("D69").Value = "LDA visualization from paper"
$ws.Range("D67").Value = "LDA2Vec: does not work (import problem)"

# Match the row height used by the neighboring TODO-only rows
$ws.Range("D67:D69").RowHeight = 17

# --- B58 ("PCA & cluster in DB, ...") loses its one-off bold/border style  ---
# and becomes identical in formatting to B57 just above it (copy format only).
$ws.Range("B57").Copy()
$ws.Range("B58").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update view state: scrolled one row further down, new active cell ---
$ws.Range("C58").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 50
$window.ScrollColumn = 1
